$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet index 1 / rId1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 2989
$ws1.Range("F5").Value = 112
$ws1.Range("F6").Value = 192
$ws1.Range("F7").Value = 1638
$ws1.Range("F8").Value = 1608
$ws1.Range("F10").Value = 354
$ws1.Range("F19").Value = 12
$ws1.Range("F20").Value = 38
$ws1.Range("F21").Value = 9
$ws1.Range("F22").Value = 355
$ws1.Range("F23").Value = 147
$ws1.Range("F24").Value = 94
$ws1.Range("F26").Value = 2000
$ws1.Range("F27").Value = 51
$ws1.Range("F28").Value = 457
$ws1.Range("F29").Value = 14
$ws1.Range("F30").Value = 169
$ws1.Range("F31").Value = 571
$ws1.Range("F32").Value = 223
$ws1.Range("F35").Value = 490
$ws1.Range("F36").Value = 5

# Sheet "全部类型" (sheet index 4 / rId4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 2989
$ws4.Range("F5").Value = 112
$ws4.Range("F6").Value = 192
$ws4.Range("F7").Value = 1638
$ws4.Range("F8").Value = 1608
$ws4.Range("F10").Value = 354
$ws4.Range("F19").Value = 12
$ws4.Range("F20").Value = 38
$ws4.Range("F21").Value = 9
$ws4.Range("F22").Value = 356
$ws4.Range("F23").Value = 147
$ws4.Range("F26").Value = 2000
$ws4.Range("F27").Value = 51
$ws4.Range("F28").Value = 457
$ws4.Range("F29").Value = 14
$ws4.Range("F30").Value = 170
$ws4.Range("F31").Value = 571
$ws4.Range("F32").Value = 223
$ws4.Range("F35").Value = 490
$ws4.Range("F36").Value = 5

Write-Host "Done updating F column values."
